$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.949.55'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '3.417.40'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  -0.08%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '409.21'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.36%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '128.73'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -2.44%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.626'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +6.58%  '
$ws.Range("E8").Value = '  -0.05%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.748'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +11.89%  '
$ws.Range("E10").Value = '  +20.65%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '42.44'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.34%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.0000220'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +71.97%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.141'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("D14").Value = '3.961.71'
$ws.Range("E14").Value = '  +1.47%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '8.91'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +6.30%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '21.02'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +6.30%  '
$ws.Range("D17").Value = '3.428.06'
$ws.Range("E17").Value = '  +1.35%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '12.31'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +12.86%  '
$ws.Range("E19").Value = '  +4.68%  '
$ws.Range("D20").Value = '61.940.61'
$ws.Range("E20").Value = '  +0.20%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '401.01'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +26.92%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '89.62'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +6.34%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '3.18'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.67%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '13.14'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +3.41%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '3.22'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +3.48%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '32.66'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +10.57%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '8.66'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +5.71%  '
$ws.Range("E28").Value = '  +0.37%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '7.59'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.71%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.70'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.62%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.118'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +2.35%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.172'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.06%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '11.85'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +4.43%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '42.97'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +2.87%  '
$ws.Range("E35").Value = '  +0.66%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.0499'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +3.78%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '53.74'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +3.77%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.15%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '3.36'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.92%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.133'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +7.30%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.91'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.77%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.312'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +6.64%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '141.09'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +1.42%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.97'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.54%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '4.09'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +3.49%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.41'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +8.89%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '16.63'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.18%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '21.69'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.58%  '
$ws.Range("D49").Value = '2.118.41'
$ws.Range("E49").Value = '  -0.29%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.35'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +2.13%  '
$ws.Range("E51").Value = '  +15.87%  '
